# "Added PDF versions to site"
#
# The site's build now generates a PDF straight from the markdown
# source, so the old "Word version of this document" link in the
# Additional Resources list is removed. Drop that whole bullet
# (paragraph + hyperlink) and leave everything else untouched.

$d = $word.ActiveDocument

$target = "Word version of this document"

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd("`r", "`a", "`n")
    if ($text -eq $target) {
        # Remove the paragraph (including its end-of-paragraph mark)
        # so the remaining bullets close up without leaving a gap.
        $p.Range.Delete()
    }
}
